$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header suffixes: "_old" -> "_FV2404", "_new" -> "_FV2410"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        $v = $v -replace '_old$', '_FV2404'
        $v = $v -replace '_new$', '_FV2410'
        $cell.Value = $v
    }
}

# 2) Freeze the header row (split at row 1, frozen top-left pane)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Add an Excel Table (ListObject) over the used range, with a header row
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$lo.Name = "Table1"
